$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New boss "unicorn" achievements & tweaks, alongside the existing "qiongqi"
# (穷奇) boss achievements. The writes below are intentionally ordered so
# that every brand-new text value first appears in the same sequence the
# authored workbook uses, since new entries get appended to the shared
# string table in first-use order.

# Row 8: new achievement "怀疑人生" (jungle exploration) replacing the
# generic templated "首次通关" row; icon -> trees.
$ws.Range("G8").Value = "trees"
$ws.Range("B8").Value = "怀疑人生"

# Row 6: description wording tweak ("连续遇到2条溪流" -> "连续遇到2次【溪流】").
$ws.Range("C6").Value = "连续遇到2次【溪流】"

$ws.Range("C8").Value = "累计遇到过8次【丛林】"

# Row 5: new achievement "快枪手" (Quick Draw) replacing "快速游戏".
$ws.Range("B5").Value = "快枪手"

# Row 4: new achievement "屠龙人" (Dragon Slayer) replacing "击杀穷奇".
$ws.Range("B4").Value = "屠龙人"

# Row 9: new achievement "另一个传说" - find & defeat the new unicorn boss,
# replacing the generic templated "首次通关" row.
$ws.Range("B9").Value = "另一个传说"
$ws.Range("C9").Value = "找到并击败另一只神兽"
$ws.Range("G9").Value = "bossunicorn"
$ws.Range("L9").Value = "star"

# Icon columns for rows 10-11.
$ws.Range("L10").Value = "lv51"
$ws.Range("L11").Value = "lv41"

# Row 6 icon changed from head3 to water2.
$ws.Range("L6").Value = "water2"

# Row 8 icon changed from hero1 to tree.
$ws.Range("L8").Value = "tree"

# Remaining edits that reuse existing shared strings / plain numbers.
$ws.Range("H9").Value = "reward"
$ws.Range("J8").Value = 8
$ws.Range("D9").Value = 5
$ws.Range("I9").Value = 1

# Selection moved to L8 (matches author's last-edited cell).
$ws.Range("L8").Select()
